# Applies the "Undisclosed Share Scheme" / Philippines-land-interest edits
# to the Form E Questionnaire document.
#
# NOTE: the source document's currency sign is mojibake - every "£" in the
# body text is actually the two characters "Â£" (0xC2 0xA3), not a plain
# pound sign.  We reproduce that exact sequence below so Find() matches the
# stored text.

$d = $word.ActiveDocument

function WholeReplace($findText, $replaceText) {
    $ok = $d.Content.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
    if (-not $ok) {
        Write-Output "WARNING: could not find [$findText]"
    }
    return $ok
}

# Finds $findText inside paragraph number $paraIndex's range and returns the
# matching sub-range (or $null if not found), without touching text outside
# that paragraph.
function Find-InParagraph($paraIndex, $findText) {
    $rng = $d.Paragraphs($paraIndex).Range.Duplicate
    $ok = $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        Write-Output "WARNING: could not find [$findText] in paragraph $paraIndex"
        return $null
    }
    return $rng
}

# =======================================================================
# Section "As to Section 2.4 & 2.15 (Undisclosed Investments / M+G
# Payment)" -> "... (Undisclosed Share Scheme / M+G Payment)"
# =======================================================================

# 1) Heading text
WholeReplace "As to Section 2.4 & 2.15 (Undisclosed Investments / M+G Payment)" `
             "As to Section 2.4 & 2.15 (Undisclosed Share Scheme / M+G Payment)"

# 2) Intro sentence: append the "aware of share/option scheme" clause
WholeReplace " on 4 June 2025. Please:" `
             " on 4 June 2025. The Respondent is aware that the Applicant participates (or participated) in an employee share/option scheme. Please:"

# 3) Item (a): bold the "a" tab-marker run, then replace the body text with
#    a bold "Provide full details" lead-in plus a normal continuation.
$paraA = $d.Paragraphs(17)
$pStart = $paraA.Range.Start
$markerRange = $d.Range($pStart, $pStart + 3)
$markerRange.Bold = 1

$body = Find-InParagraph 17 "Explain the nature of this payment (dividend, share sale, pension drawdown?)."
if ($body -ne $null) {
    $body.Text = "Provide full details of any Share Incentive Plan (SIP), SAYE, or Stock Option scheme operated by your current or former employer, including the number of units held (vested or unvested)."
    $boldRange = $d.Range($body.Start, $body.Start + 20)
    $boldRange.Bold = 1
}

# 4) Item (b): M&G continuing interest -> confirm if payment was a dividend
$body = Find-InParagraph 18 "Confirm if the Applicant holds any continuing interest in an M&G investment, pension, or employee share scheme."
if ($body -ne $null) {
    $body.Text = 'Confirm if the Â£1,912.47 M+G payment was a dividend or payout from this scheme.'
}

# 5) Item (c): "Provide documentation" -> "Provide the scheme statement";
#    trailing clause rewritten.
$body = Find-InParagraph 19 "Provide documentation"
if ($body -ne $null) {
    $body.Text = "Provide the scheme statement"
}
$body = Find-InParagraph 19 " evidencing the source of this payment (e.g., dividend voucher, contract note)."
if ($body -ne $null) {
    $body.Text = " or documentation showing current holdings and value."
}

# =======================================================================
# Section "As to Section 2.10 (Property & Assets Outside Jurisdiction)"
# =======================================================================

# 6) Intro sentence rewrite
$body = Find-InParagraph 27 'The Applicant states "Nil" foreign assets. Bank statements show remittances to the Philippines totaling Â£3,524.79 (Aug 24 - Oct 25) via Remitly, Ms. Joylyn Gray, and Ms. Lucy Timog. Please:'
if ($body -ne $null) {
    $body.Text = 'The Applicant states "Nil" foreign assets. However, bank statements show remittances to the Philippines totaling Â£3,524.79 (Aug 24 - Oct 25). The Respondent is aware of the Applicant''s prior stated intention to purchase land. Please:'
}

# 7) Item (a): lead text -> "Specifically ", bold text -> "confirm or deny",
#    tail -> long sentence about legal/beneficial interest.
$body = Find-InParagraph 28 "Confirm if any funds relate to the purchase or improvement of land/property. If so, "
if ($body -ne $null) {
    $body.Text = "Specifically "
}
$body = Find-InParagraph 28 "provide deeds and valuations"
if ($body -ne $null) {
    $body.Text = "confirm or deny"
}
$body = Find-InParagraph 28 "."
if ($body -ne $null) {
    $body.Text = " whether the Applicant holds any legal or beneficial interest (solely or jointly) in land or property in the Philippines."
}

# 8) Item (b): family-remittance-affordability text -> 'If "Yes", ' plus a
#    freshly (re)created bold "provide deeds and valuations" + "." tail -
#    i.e. what used to live in item (a) before the rewrite above.
$body = Find-InParagraph 29 "If these are family remittances, explain the affordability of"
if ($body -ne $null) {
    $pEnd = $d.Paragraphs(29).Range.End
    $full = $d.Range($body.Start, $pEnd - 1)
    $full.Text = 'If "Yes", '

    $pEnd2 = $d.Paragraphs(29).Range.End
    $insertPoint = $d.Range($pEnd2 - 1, $pEnd2 - 1)
    $insertPoint.InsertAfter("provide deeds and valuations")
    $boldStart = $pEnd2 - 1
    $boldEnd = $boldStart + 29
    $boldRange = $d.Range($boldStart, $boldEnd)
    $boldRange.Bold = 1
    $dotPoint = $d.Range($boldEnd, $boldEnd)
    $dotPoint.InsertAfter(".")
}

# 9) Item (c): merge the bold "Provide transfer receipts" run and the
#    normal tail run into a single normal-weight sentence.
$body = Find-InParagraph 30 "Provide transfer receipts"
if ($body -ne $null) {
    $body.Bold = 0
    $body.Text = 'If "No", explain the specific purpose of the Â£3,524 remittances and provide transfer receipts showing the ultimate recipient.'
}
$body = Find-InParagraph 30 " showing the ultimate recipient of these funds."
if ($body -ne $null) {
    $body.Text = ""
}

Write-Output "done"
